# Insert a new row at position 618 (shifts existing rows 618-697 down to 619-698)
# and populate it with the new "Apio" price record for 2023-08-16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(618).Insert()

$ws.Cells.Item(618, 1).Value  = 3
$ws.Cells.Item(618, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(618, 3).Value  = "Coquimbo"
$ws.Cells.Item(618, 4).Value  = 45154
$ws.Cells.Item(618, 5).Value  = 5
$ws.Cells.Item(618, 6).Value  = 100112017
$ws.Cells.Item(618, 7).Value  = "Apio"
$ws.Cells.Item(618, 8).Value  = "Americana (o)"
$ws.Cells.Item(618, 9).Value  = "Primera"
$ws.Cells.Item(618, 10).Value = 220
$ws.Cells.Item(618, 11).Value = 8000
$ws.Cells.Item(618, 12).Value = 8500
$ws.Cells.Item(618, 13).Value = 8227
$ws.Cells.Item(618, 14).Value = "$/docena de matas"
$ws.Cells.Item(618, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(618, 16).Value = 1371
$ws.Cells.Item(618, 17).Value = 6
$ws.Cells.Item(618, 18).Value = "Hortaliza"
